# Work on the workbook / active sheet (Sheet1) exactly as the user had it open.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data entered by the user:
#   B1 -> a plain number
#   B3 -> a new text value ("rrrrrr"), which becomes a new shared-string entry
$ws.Range("B1").Value = 11221223
$ws.Range("B3").Value = "rrrrrr"

# Let Excel size column B to fit its new contents (mirrors the bestFit/customWidth
# behaviour recorded for column B in the saved workbook).
$ws.Columns.Item(2).AutoFit() | Out-Null

# The user finished by leaving the selection on the newly entered cell B3.
$ws.Range("B3").Select() | Out-Null
